$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp cell (A1) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 19 de Agosto de 2020 a las 21:49"

# --- Re-rank country rows whose order changed (string swaps) + refresh daily COVID figures ---
# Row 4: data refresh
$ws.Cells.Item(4, 2).Value = 5681160
$ws.Cells.Item(4, 3).Value = 25186
$ws.Cells.Item(4, 4).Value = 3036111
$ws.Cells.Item(4, 5).Value = 2469235
$ws.Cells.Item(4, 7).Value = 740
$ws.Cells.Item(4, 8).Value = 175814

# Row 6: data refresh
$ws.Cells.Item(6, 2).Value = 2835822
$ws.Cells.Item(6, 3).Value = 69196
$ws.Cells.Item(6, 4).Value = 2096068
$ws.Cells.Item(6, 5).Value = 685760
$ws.Cells.Item(6, 7).Value = 980
$ws.Cells.Item(6, 8).Value = 53994

# Row 22: data refresh
$ws.Cells.Item(22, 2).Value = 229459
$ws.Cells.Item(22, 3).Value = 1354
$ws.Cells.Item(22, 5).Value = 16245
$ws.Cells.Item(22, 7).Value = 9
$ws.Cells.Item(22, 8).Value = 9314

# Row 23: data refresh
$ws.Cells.Item(23, 2).Value = 225043
$ws.Cells.Item(23, 3).Value = 3776
$ws.Cells.Item(23, 5).Value = 110510
$ws.Cells.Item(23, 7).Value = 17
$ws.Cells.Item(23, 8).Value = 30468

# Row 29: now "Ecuador"
$ws.Cells.Item(29, 1).Value = "Ecuador"
$ws.Cells.Item(29, 2).Value = 104475
$ws.Cells.Item(29, 3).Value = 1534
$ws.Cells.Item(29, 4).Value = 87277
$ws.Cells.Item(29, 5).Value = 11052
$ws.Cells.Item(29, 7).Value = 41
$ws.Cells.Item(29, 8).Value = 6146

# Row 30: now "Kazajistan"
$ws.Cells.Item(30, 1).Value = "Kazajistan"
$ws.Cells.Item(30, 2).Value = 103571
$ws.Cells.Item(30, 3).Value = 271
$ws.Cells.Item(30, 4).Value = 86286
$ws.Cells.Item(30, 5).Value = 15870
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 1415

# Row 31: now "Bolivia"
$ws.Cells.Item(31, 1).Value = "Bolivia"
$ws.Cells.Item(31, 2).Value = 103019
$ws.Cells.Item(31, 3).Value = 1796
$ws.Cells.Item(31, 4).Value = 38623
$ws.Cells.Item(31, 5).Value = 60224
$ws.Cells.Item(31, 7).Value = 49
$ws.Cells.Item(31, 8).Value = 4172

# Row 32: data refresh
$ws.Cells.Item(32, 2).Value = 97783
$ws.Cells.Item(32, 3).Value = 1374
$ws.Cells.Item(32, 4).Value = 73084
$ws.Cells.Item(32, 5).Value = 23918
$ws.Cells.Item(32, 7).Value = 73
$ws.Cells.Item(32, 8).Value = 781

# Row 55: data refresh
$ws.Cells.Item(55, 2).Value = 43094
$ws.Cells.Item(55, 3).Value = 101
$ws.Cells.Item(55, 4).Value = 40963
$ws.Cells.Item(55, 5).Value = 1875
$ws.Cells.Item(55, 7).Value = 8
$ws.Cells.Item(55, 8).Value = 256

# Row 64: data refresh
$ws.Cells.Item(64, 2).Value = 34058
$ws.Cells.Item(64, 3).Value = 1336
$ws.Cells.Item(64, 4).Value = 13308
$ws.Cells.Item(64, 5).Value = 20150
$ws.Cells.Item(64, 7).Value = 28
$ws.Cells.Item(64, 8).Value = 600

# Row 67: now "Costa Rica"
$ws.Cells.Item(67, 1).Value = "Costa Rica"
$ws.Cells.Item(67, 2).Value = 30409
$ws.Cells.Item(67, 3).Value = 766
$ws.Cells.Item(67, 4).Value = 9660
$ws.Cells.Item(67, 5).Value = 20428
$ws.Cells.Item(67, 7).Value = 7
$ws.Cells.Item(67, 8).Value = 321

# Row 68: now "Serbia"
$ws.Cells.Item(68, 1).Value = "Serbia"
$ws.Cells.Item(68, 2).Value = 30048
$ws.Cells.Item(68, 3).Value = 158
$ws.Cells.Item(68, 4).Value = 27702
$ws.Cells.Item(68, 5).Value = 1662
$ws.Cells.Item(68, 7).Value = 3
$ws.Cells.Item(68, 8).Value = 684

# Row 75: data refresh
$ws.Cells.Item(75, 2).Value = 18624
$ws.Cells.Item(75, 3).Value = 25
$ws.Cells.Item(75, 5).Value = 1678

# Row 101: data refresh
$ws.Cells.Item(101, 2).Value = 7566
$ws.Cells.Item(101, 3).Value = 67
$ws.Cells.Item(101, 4).Value = 6813
$ws.Cells.Item(101, 5).Value = 629

# Row 103: data refresh
$ws.Cells.Item(103, 2).Value = 6829
$ws.Cells.Item(103, 3).Value = 40
$ws.Cells.Item(103, 4).Value = 6094
$ws.Cells.Item(103, 5).Value = 577
$ws.Cells.Item(103, 7).Value = 1
$ws.Cells.Item(103, 8).Value = 158

# Row 105: data refresh
$ws.Cells.Item(105, 2).Value = 5643
$ws.Cells.Item(105, 3).Value = 265
$ws.Cells.Item(105, 4).Value = 4442
$ws.Cells.Item(105, 5).Value = 1051
$ws.Cells.Item(105, 7).Value = 9
$ws.Cells.Item(105, 8).Value = 150

# Row 107: data refresh
$ws.Cells.Item(107, 2).Value = 5240
$ws.Cells.Item(107, 3).Value = 47
$ws.Cells.Item(107, 4).Value = 2857
$ws.Cells.Item(107, 5).Value = 2219
$ws.Cells.Item(107, 7).Value = 1
$ws.Cells.Item(107, 8).Value = 164

# Row 115: data refresh
$ws.Cells.Item(115, 2).Value = 4058
$ws.Cells.Item(115, 3).Value = 69
$ws.Cells.Item(115, 4).Value = 2611
$ws.Cells.Item(115, 5).Value = 1368
$ws.Cells.Item(115, 7).Value = 3
$ws.Cells.Item(115, 8).Value = 79

# Row 116: data refresh
$ws.Cells.Item(116, 2).Value = 3850
$ws.Cells.Item(116, 3).Value = 19
$ws.Cells.Item(116, 5).Value = 2148
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = 77

# Row 117: data refresh
$ws.Cells.Item(117, 2).Value = 3482
$ws.Cells.Item(117, 3).Value = 74
$ws.Cells.Item(117, 4).Value = 2863
$ws.Cells.Item(117, 5).Value = 531

# Row 126: data refresh
$ws.Cells.Item(126, 2).Value = 2667
$ws.Cells.Item(126, 3).Value = 1
$ws.Cells.Item(126, 4).Value = 1993
$ws.Cells.Item(126, 5).Value = 549

# Row 131: data refresh
$ws.Cells.Item(131, 2).Value = 2427
$ws.Cells.Item(131, 3).Value = 113
$ws.Cells.Item(131, 4).Value = 1395
$ws.Cells.Item(131, 5).Value = 972
$ws.Cells.Item(131, 7).Value = 3
$ws.Cells.Item(131, 8).Value = 60

# Row 132: now "Gambia"
$ws.Cells.Item(132, 1).Value = "Gambia"
$ws.Cells.Item(132, 2).Value = 2288
$ws.Cells.Item(132, 3).Value = 172
$ws.Cells.Item(132, 4).Value = 435
$ws.Cells.Item(132, 5).Value = 1776
$ws.Cells.Item(132, 7).Value = 14
$ws.Cells.Item(132, 8).Value = 77

# Row 133: now "Estonia"
$ws.Cells.Item(133, 1).Value = "Estonia"
$ws.Cells.Item(133, 2).Value = 2207
$ws.Cells.Item(133, 3).Value = 7
$ws.Cells.Item(133, 4).Value = 2002
$ws.Cells.Item(133, 5).Value = 142
$ws.Cells.Item(133, 8).Value = 63

# Row 134: now "Guinea-Bisau"
$ws.Cells.Item(134, 1).Value = "Guinea-Bisau"
$ws.Cells.Item(134, 2).Value = 2117
$ws.Cells.Item(134, 4).Value = 1015
$ws.Cells.Item(134, 5).Value = 1069
$ws.Cells.Item(134, 8).Value = 33

# Row 135: data refresh
$ws.Cells.Item(135, 2).Value = 2095
$ws.Cells.Item(135, 3).Value = 32
$ws.Cells.Item(135, 5).Value = 366

# Row 139: now "Siria"
$ws.Cells.Item(139, 1).Value = "Siria"
$ws.Cells.Item(139, 2).Value = 1927
$ws.Cells.Item(139, 3).Value = 83
$ws.Cells.Item(139, 4).Value = 445
$ws.Cells.Item(139, 5).Value = 1404
$ws.Cells.Item(139, 7).Value = 5
$ws.Cells.Item(139, 8).Value = 78

# Row 140: now "Yemen"
$ws.Cells.Item(140, 1).Value = "Yemen"
$ws.Cells.Item(140, 2).Value = 1892
$ws.Cells.Item(140, 3).Value = 3
$ws.Cells.Item(140, 4).Value = 1055
$ws.Cells.Item(140, 5).Value = 298
$ws.Cells.Item(140, 7).Value = 2
$ws.Cells.Item(140, 8).Value = 539

# Row 151: now "Aruba"
$ws.Cells.Item(151, 1).Value = "Aruba"
$ws.Cells.Item(151, 2).Value = 1296
$ws.Cells.Item(151, 3).Value = 91
$ws.Cells.Item(151, 4).Value = 277
$ws.Cells.Item(151, 5).Value = 1014
$ws.Cells.Item(151, 7).Value = 1
$ws.Cells.Item(151, 8).Value = 5

# Row 152: now "Burkina Faso"
$ws.Cells.Item(152, 1).Value = "Burkina Faso"
$ws.Cells.Item(152, 2).Value = 1285
$ws.Cells.Item(152, 3).Value = 5
$ws.Cells.Item(152, 4).Value = 1023
$ws.Cells.Item(152, 5).Value = 207
$ws.Cells.Item(152, 8).Value = 55

# Row 153: now "Liberia"
$ws.Cells.Item(153, 1).Value = "Liberia"
$ws.Cells.Item(153, 2).Value = 1282
$ws.Cells.Item(153, 4).Value = 803
$ws.Cells.Item(153, 5).Value = 397
$ws.Cells.Item(153, 8).Value = 82

# Row 159: data refresh
$ws.Cells.Item(159, 2).Value = 971
$ws.Cells.Item(159, 3).Value = 1
$ws.Cells.Item(159, 4).Value = 868
$ws.Cells.Item(159, 5).Value = 27

# Row 171: data refresh
$ws.Cells.Item(171, 2).Value = 422
$ws.Cells.Item(171, 3).Value = 6
$ws.Cells.Item(171, 5).Value = 85

# Row 192: data refresh
$ws.Cells.Item(192, 2).Value = 132
$ws.Cells.Item(192, 3).Value = 5
$ws.Cells.Item(192, 5).Value = 6

# Row 213: now "Montserrat"
$ws.Cells.Item(213, 1).Value = "Montserrat"
$ws.Cells.Item(213, 4).Value = 12
$ws.Cells.Item(213, 8).Value = 1

# Row 214: now "Islas Malvinas"
$ws.Cells.Item(214, 1).Value = "Islas Malvinas"
$ws.Cells.Item(214, 4).Value = 13
$ws.Cells.Item(214, 8).Value = 0
